$wb = $excel.ActiveWorkbook

$wsOptions = $wb.Worksheets.Item("Geometry options")
$wsCalc = $wb.Worksheets.Item("geometry calculation")

# --- Core data edit -----------------------------------------------------
# BC-LEEP archetype (row 10) on "Geometry options": set the basement wall
# height (column Q) from 0 to 2.4384 m (8 ft). Every other changed cell in
# the workbook is a formula that depends on this value (directly or via
# the named range Opt_Bsm_Height), so Excel's recalculation engine
# propagates the change automatically.
$wsOptions.Range("Q10").Value = 2.4384

# --- Cosmetic view-state changes -----------------------------------------
# Column C on "geometry calculation" was widened slightly so the larger
# recalculated numbers still display comfortably.
$wsCalc.Columns.Item(3).ColumnWidth = 10.59

# The user ended up with "geometry calculation" as the active sheet/tab,
# scrolled down a bit with C4 selected; "Geometry options" was left with
# Q10 (the cell just edited) selected but no longer the active tab.
[void]$wsOptions.Range("Q10").Select()
[void]$wsCalc.Activate()
[void]$wsCalc.Range("C4").Select()
